$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136, shifting existing rows 136-160 down to 137-161
$ws.Rows.Item(136).Insert()

# Populate the new row 136 with the new record (copy static columns from what is now row 137,
# then overwrite the values that differ per the diff)
$ws.Range("A136").Value = 10
$ws.Range("B136").Value = "Vega Modelo de Temuco"
$ws.Range("C136").Value = "La Araucanía"
$ws.Range("D136").Value = 45275
$ws.Range("D136").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E136").Value = 9
$ws.Range("F136").Value = "Fruta"
$ws.Range("G136").Value = 100101
$ws.Range("H136").Value = "Berries"
$ws.Range("I136").Value = 100101001
$ws.Range("J136").Value = "Arándano (blue)"
$ws.Range("K136").Value = "Sin especificar"
$ws.Range("L136").Value = "Primera"
$ws.Range("M136").Value = 380
$ws.Range("N136").Value = 2700
$ws.Range("O136").Value = 2800
$ws.Range("P136").Value = 2753
$ws.Range("Q136").Value = "$/kilo"
$ws.Range("R136").Value = "Región del Maule"
$ws.Range("S136").Value = 2753
$ws.Range("T136").Value = 1
